$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 14:12"

# Update province/city data rows that shifted or changed due to refreshed data
$ws.Cells.Item(16, 1).Value = "Malaga"
$ws.Cells.Item(16, 2).Value = 1006
$ws.Cells.Item(16, 3).Value = 61
$ws.Cells.Item(16, 4).Value = 899
$ws.Cells.Item(16, 5).Value = 46

$ws.Cells.Item(17, 1).Value = "Toledo"
$ws.Cells.Item(17, 2).Value = 965
$ws.Cells.Item(17, 3).Value = 95
$ws.Cells.Item(17, 4).Value = 860
$ws.Cells.Item(17, 5).Value = 78

$ws.Cells.Item(18, 1).Value = "Aragon"
$ws.Cells.Item(18, 2).Value = 907
$ws.Cells.Item(18, 3).Value = 29
$ws.Cells.Item(18, 4).Value = 838
$ws.Cells.Item(18, 5).Value = 40

$ws.Cells.Item(25, 1).Value = "Granada"
$ws.Cells.Item(25, 2).Value = 711
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 676
$ws.Cells.Item(25, 5).Value = 34

$ws.Cells.Item(26, 1).Value = "Sevilla"
$ws.Cells.Item(26, 2).Value = 708
$ws.Cells.Item(26, 3).Value = 8
$ws.Cells.Item(26, 4).Value = 675
$ws.Cells.Item(26, 5).Value = 25

$ws.Cells.Item(27, 1).Value = "Murcia"
$ws.Cells.Item(27, 2).Value = 687
$ws.Cells.Item(27, 3).Value = 12
$ws.Cells.Item(27, 4).Value = 660
$ws.Cells.Item(27, 5).Value = 15

$ws.Cells.Item(28, 1).Value = "Albacete"
$ws.Cells.Item(28, 2).Value = 666
$ws.Cells.Item(28, 3).Value = 95
$ws.Cells.Item(28, 4).Value = 592
$ws.Cells.Item(28, 5).Value = 66

$ws.Cells.Item(29, 1).Value = "Salamanca"
$ws.Cells.Item(29, 2).Value = 629
$ws.Cells.Item(29, 3).Value = 73
$ws.Cells.Item(29, 4).Value = 502
$ws.Cells.Item(29, 5).Value = 54

$ws.Cells.Item(35, 1).Value = "Jaen"
$ws.Cells.Item(35, 2).Value = 414
$ws.Cells.Item(35, 3).Value = 7
$ws.Cells.Item(35, 4).Value = 392
$ws.Cells.Item(35, 5).Value = 15

$ws.Cells.Item(36, 1).Value = "Castello/Castellon"
$ws.Cells.Item(36, 2).Value = 412
$ws.Cells.Item(36, 3).Value = 4
$ws.Cells.Item(36, 4).Value = 389
$ws.Cells.Item(36, 5).Value = 19

$ws.Cells.Item(37, 1).Value = "Badajoz"
$ws.Cells.Item(37, 2).Value = 390
$ws.Cells.Item(37, 3).Value = 29
$ws.Cells.Item(37, 4).Value = 352
$ws.Cells.Item(37, 5).Value = 9

$ws.Cells.Item(38, 1).Value = "Segovia"
$ws.Cells.Item(38, 2).Value = 361
$ws.Cells.Item(38, 3).Value = 62
$ws.Cells.Item(38, 4).Value = 262
$ws.Cells.Item(38, 5).Value = 37

$ws.Cells.Item(39, 1).Value = "Cordoba"
$ws.Cells.Item(39, 2).Value = 359
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = 350
$ws.Cells.Item(39, 5).Value = 8

$ws.Cells.Item(40, 1).Value = "Cadiz"
$ws.Cells.Item(40, 2).Value = 348
$ws.Cells.Item(40, 3).Value = 6
$ws.Cells.Item(40, 4).Value = 336
$ws.Cells.Item(40, 5).Value = 6

$ws.Cells.Item(41, 1).Value = "Soria"
$ws.Cells.Item(41, 2).Value = 339
$ws.Cells.Item(41, 3).Value = 32
$ws.Cells.Item(41, 4).Value = 291
$ws.Cells.Item(41, 5).Value = 16

$ws.Cells.Item(42, 1).Value = "Ourense"
$ws.Cells.Item(42, 2).Value = 321
$ws.Cells.Item(42, 3).Value = 67
$ws.Cells.Item(42, 4).Value = 302
$ws.Cells.Item(42, 5).Value = 5

$ws.Cells.Item(49, 1).Value = "Almeria"
$ws.Cells.Item(49, 2).Value = 142
$ws.Cells.Item(49, 3).Value = 6
$ws.Cells.Item(49, 4).Value = 128
$ws.Cells.Item(49, 5).Value = 8

$ws.Cells.Item(50, 1).Value = "Palencia"
$ws.Cells.Item(50, 2).Value = 139
$ws.Cells.Item(50, 3).Value = 14
$ws.Cells.Item(50, 4).Value = 120
$ws.Cells.Item(50, 5).Value = 5

$ws.Cells.Item(52, 1).Value = "Huelva"
$ws.Cells.Item(52, 2).Value = 105
$ws.Cells.Item(52, 3).Value = 2
$ws.Cells.Item(52, 4).Value = 101
$ws.Cells.Item(52, 5).Value = 2

$ws.Cells.Item(54, 1).Value = "Melilla"
$ws.Cells.Item(54, 2).Value = 42
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(54, 4).Value = 41
$ws.Cells.Item(54, 5).Value = 1

$ws.Cells.Item(59, 1).Value = "Ceuta"
$ws.Cells.Item(59, 2).Value = 17
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 16
$ws.Cells.Item(59, 5).Value = 1
